# Added SETU document types for TICC-319
#
# Appends six rows (254-259) to the bottom of the "Document Type" sheet:
# SETU HR-XML StaffingOrder / StaffingOrder Status / HumanResource /
# HumanResource Status / Assignment / Assignment Status, each v1.4 (or
# v1.4.1 for Assignment), release "8.9", state "active", comment
# "TICC-319", domain community "Extended use". This mirrors the existing
# "SETU HR-XML Timecard" block added for the previous ticket TICC-266
# (row 220), which is reused below purely as a formatting donor.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Document Type")

$templateRow = 220

# (row, @{col = value}) pairs, listed in the same left-to-right /
# row-by-row order the values were actually typed in, so brand-new
# shared-string entries land in the same sequence as the authored sheet.
$firstNewRow = 254
$rows = 254..259

$values = @{
    254 = @{ A = "SETU HR-XML StaffingOrder v1.4" }
    255 = @{ A = "SETU HR-XML StaffingOrder Status v1.4" }
    256 = @{ A = "SETU HR-XML HumanResource v1.4" }
    257 = @{ A = "SETU HR-XML HumanResource Status v1.4" }
    258 = @{ A = "SETU HR-XML Assignment v1.4.1" }
    259 = @{ A = "SETU HR-XML Assignment Status v1.4" }
}

$cValues = @{
    254 = "http://ns.hr-xml.org/2007-04-15::StaffingOrder##hr-xml@nl-1.4::2.5"
    255 = "http://ns.hr-xml.org/2007-04-15::StaffingOrder##hr-xml:status@nl-1.4::2.5"
    256 = "http://ns.hr-xml.org/2007-04-15::HumanResource##hr-xml@nl-1.4::2.5"
    257 = "http://ns.hr-xml.org/2007-04-15::HumanResource##hr-xml:status@nl-1.4::2.5"
    258 = "http://ns.hr-xml.org/2007-04-15::Assignment##hr-xml@nl-1.4.1::2.5"
    259 = "http://ns.hr-xml.org/2007-04-15::Assignment##hr-xml:status@nl-1.4::2.5"
}

$mValues = @{
    254 = "cenbii-procid-ubl::urn:fdc:hr-xml:2007:staffingorder:1.0"
    255 = "cenbii-procid-ubl::urn:fdc:hr-xml:2007:staffingorder:1.0"
    256 = "cenbii-procid-ubl::urn:fdc:hr-xml:2007:humanresource:1.0"
    257 = "cenbii-procid-ubl::urn:fdc:hr-xml:2007:humanresource:1.0"
    258 = "cenbii-procid-ubl::urn:fdc:hr-xml:2007:assignment:1.0"
    259 = "cenbii-procid-ubl::urn:fdc:hr-xml:2007:assignment:1.0"
}

# 1) Stamp every new row with the previous SETU entry's formatting
#    (A/B/D/E/H/I/J/L/M columns only -- C/F/G/K are left untouched, same as
#    the template row, which has no cells in those columns either).
$formattedCols = @("A", "B", "D", "E", "H", "I", "J", "L", "M")
foreach ($r in $rows) {
    foreach ($col in $formattedCols) {
        $ws.Range("$col$templateRow").Copy() | Out-Null
        $ws.Range("$col$r").PasteSpecial(-4122) | Out-Null
    }
}
$excel.CutCopyMode = 0

# 2) Row 254 (StaffingOrder) filled in first, left to right.
$ws.Range("A254").Value = $values[254].A
$ws.Range("B254").Value = "busdox-docid-qns"
$ws.Range("D254").Value = "8.9"
$ws.Range("E254").Value = "active"
$ws.Range("M254").Value = $mValues[254]
$ws.Range("H254").Value = "TICC-319"
$ws.Range("I254").Value = $false
$ws.Range("J254").Value = $false
$ws.Range("L254").Value = "Extended use"
$ws.Range("C254").Value = $cValues[254]

# 3) StaffingOrder Status row -- just the profile name for now.
$ws.Range("A255").Value = $values[255].A
$ws.Range("B255").Value = "busdox-docid-qns"
$ws.Range("D255").Value = "8.9"
$ws.Range("E255").Value = "active"
$ws.Range("H255").Value = "TICC-319"
$ws.Range("I255").Value = $false
$ws.Range("J255").Value = $false
$ws.Range("L255").Value = "Extended use"

# 4) HumanResource row (name + identifier value).
$ws.Range("A256").Value = $values[256].A
$ws.Range("B256").Value = "busdox-docid-qns"
$ws.Range("C256").Value = $cValues[256]
$ws.Range("D256").Value = "8.9"
$ws.Range("E256").Value = "active"
$ws.Range("H256").Value = "TICC-319"
$ws.Range("I256").Value = $false
$ws.Range("J256").Value = $false
$ws.Range("L256").Value = "Extended use"

# 5) HumanResource Status row -- profile name only.
$ws.Range("A257").Value = $values[257].A
$ws.Range("B257").Value = "busdox-docid-qns"
$ws.Range("D257").Value = "8.9"
$ws.Range("E257").Value = "active"
$ws.Range("H257").Value = "TICC-319"
$ws.Range("I257").Value = $false
$ws.Range("J257").Value = $false
$ws.Range("L257").Value = "Extended use"

# 6) Assignment row (name + identifier value).
$ws.Range("A258").Value = $values[258].A
$ws.Range("B258").Value = "busdox-docid-qns"
$ws.Range("C258").Value = $cValues[258]
$ws.Range("D258").Value = "8.9"
$ws.Range("E258").Value = "active"
$ws.Range("H258").Value = "TICC-319"
$ws.Range("I258").Value = $false
$ws.Range("J258").Value = $false
$ws.Range("L258").Value = "Extended use"

# 7) Assignment Status row -- profile name only.
$ws.Range("A259").Value = $values[259].A
$ws.Range("B259").Value = "busdox-docid-qns"
$ws.Range("D259").Value = "8.9"
$ws.Range("E259").Value = "active"
$ws.Range("H259").Value = "TICC-319"
$ws.Range("I259").Value = $false
$ws.Range("J259").Value = $false
$ws.Range("L259").Value = "Extended use"

# 8) Go back and fill in the "Associated Process/Profile Identifier(s)"
#    column for the two still-missing parent rows...
$ws.Range("M256").Value = $mValues[256]
$ws.Range("M258").Value = $mValues[258]

# ...then the "Status" rows re-use the already-minted M254/M256/M258
# strings.
$ws.Range("M255").Value = $mValues[255]
$ws.Range("M257").Value = $mValues[257]
$ws.Range("M259").Value = $mValues[259]

# 9) Finally, fill the identifier values for the three "Status" rows.
$ws.Range("C255").Value = $cValues[255]
$ws.Range("C257").Value = $cValues[257]
$ws.Range("C259").Value = $cValues[259]

# Match the author's final selection.
$ws.Range("A258").Select() | Out-Null
